# BIS-1002: removed "Internal Assignment" column from export. Expanded and fixed tests
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Internal Assignment" column contents (header + TRUE/FALSE values)
$ws.Range("O4:O8").ClearContents()

# Fix test data value in A5 ($NAME -> $$NAME)
$ws.Range("A5").Value = '$$NAME'

# Move the active selection to K14 (matches post-edit cursor position)
$ws.Range("K14").Select()
